$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 98, pushing existing rows 98:204 down to 99:205
$ws.Rows.Item(98).Insert()

# Fill in the boilerplate columns (same constant values used throughout the sheet)
$ws.Cells.Item(98, 1).Value = 11
$ws.Cells.Item(98, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(98, 3).Value = "Bíobío"
$ws.Cells.Item(98, 4).Value = 44740
$ws.Cells.Item(98, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(98, 5).Value = 8
$ws.Cells.Item(98, 6).Value = 100112040
$ws.Cells.Item(98, 7).Value = "Cilantro"
$ws.Cells.Item(98, 8).Value = "Sin especificar"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 100
$ws.Cells.Item(98, 11).Value = 11000
$ws.Cells.Item(98, 12).Value = 12000
$ws.Cells.Item(98, 13).Value = 11500
$ws.Cells.Item(98, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(98, 15).Value = "Región Metropolitana"
$ws.Cells.Item(98, 16).Value = 319
$ws.Cells.Item(98, 17).Value = 36
$ws.Cells.Item(98, 18).Value = "Hortaliza"
